$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.026.63'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '3.517.43'
$ws.Range("E3").Value = '  -1.46%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.41'
$ws.Range("E5").Value = '  -1.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.84'
$ws.Range("E6").Value = '  -1.36%  '
$ws.Range("D7").Value = '3.516.29'
$ws.Range("E7").Value = '  -1.46%  '
$ws.Range("E9").Value = '  -1.06%  '
$ws.Range("E10").Value = '  +1.05%  '
$ws.Range("E11").Value = '  +2.92%  '
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").Value = '4.115.34'
$ws.Range("E13").Value = '  -1.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.67'
$ws.Range("E14").Value = '  +2.13%  '
$ws.Range("E15").Value = '  -0.84%  '
$ws.Range("E16").Value = '  +0.72%  '
$ws.Range("D17").Value = '3.516.40'
$ws.Range("E17").Value = '  -1.57%  '
$ws.Range("D18").Value = '65.011.25'
$ws.Range("E18").Value = '  +0.62%  '
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.32'
$ws.Range("E20").Value = '  -0.50%  '
$ws.Range("E21").Value = '  -3.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '391.73'
$ws.Range("E22").Value = '  +0.69%  '
$ws.Range("E23").Value = '  -0.63%  '
$ws.Range("D24").Value = '3.659.12'
$ws.Range("E24").Value = '  -1.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.64'
$ws.Range("E25").Value = '  +0.64%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E27").Value = '  -3.73%  '
$ws.Range("E28").Value = '  +9.45%  '
$ws.Range("E29").Value = '  -1.67%  '
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("E31").Value = '  -1.24%  '
$ws.Range("E32").Value = '  -1.87%  '
$ws.Range("D33").Value = '3.523.70'
$ws.Range("E33").Value = '  -1.33%  '
$ws.Range("E34").Value = '  +0.21%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  +0.39%  '
$ws.Range("E37").Value = '  +4.66%  '
$ws.Range("E38").Value = '  +1.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").Value = '169.28'
$ws.Range("E39").Value = '  -0.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").Value = '6.92'
$ws.Range("E40").Value = '  -0.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0804'
$ws.Range("E41").Value = '  -0.73%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.821'
$ws.Range("E42").Value = '  -0.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = '26.07'
$ws.Range("E43").Value = '  -4.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("B44").Value = 'ONDO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D44").Value = '1.25'
$ws.Range("E44").Value = '  +3.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.95'
$ws.Range("E45").Value = '  +0.55%  '
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("E47").Value = '  -0.97%  '
$ws.Range("E48").Value = '  +0.77%  '
$ws.Range("E49").Value = '  -0.47%  '
$ws.Range("D50").Value = '2.437.02'
$ws.Range("E50").Value = '  -1.98%  '
$ws.Range("E51").Value = '  +3.32%  '
